# Atualização de bases das ligas, do dia: 11-06-2024 às 21:19
#
# The underlying match data for a few rows got reshuffled (rows swapped /
# rotated) while the leading index column (A) — which is just a running
# row counter — stays put. This script reads each affected row's data
# (columns B..AD) and writes it back into the row(s) it now belongs to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($ws, $row) {
    $data = @{}
    for ($c = 2; $c -le 30; $c++) {
        $data[$c] = $ws.Cells.Item($row, $c).Value2()
    }
    return $data
}

function Set-RowData($ws, $row, $data) {
    for ($c = 2; $c -le 30; $c++) {
        $ws.Cells.Item($row, $c).Value = $data[$c]
    }
}

# Rows 107 and 108 (Libertad Gran Mamore FC vs Royal Pari FC / Universitario
# De Vinto vs Bolivar) swap places with each other.
$row107 = Get-RowData $ws 107
$row108 = Get-RowData $ws 108

Set-RowData $ws 107 $row108
Set-RowData $ws 108 $row107

# Rows 148, 149 and 150 rotate: 148 -> 150, 149 -> 148, 150 -> 149
# (i.e. each row takes on the data previously held by the row above it,
# with row 148 receiving what used to be in row 150).
$row148 = Get-RowData $ws 148
$row149 = Get-RowData $ws 149
$row150 = Get-RowData $ws 150

Set-RowData $ws 148 $row149
Set-RowData $ws 149 $row150
Set-RowData $ws 150 $row148
